# The sheet currently has an extra "header-only" row (row 6) containing only
# the label "grandes regiões e unidades da federação" with no data values -
# the actual region data starts on row 7 ("norte") and continues through
# row 38 ("distrito federal").
#
# The fix removes that redundant label row entirely. Deleting it shifts every
# row below up by one: row 7 ("norte", with its numeric data) becomes the new
# row 6, row 8 becomes row 7, ..., and the former row 38 ("distrito federal")
# becomes the new (and final) row 37. The now-unused shared string
# "grandes regiões e unidades da federação" is dropped automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
